$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.496.38"
$ws.Range("E2").Value = "  +0.29%  "
$ws.Range("D3").Value = "1.939.49"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.609"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "56.97"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.359"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0851"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.32%  "
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "2.223.95"
$ws.Range("E12").Value = "  -1.86%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -6.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.806"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -6.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -6.20%  "
$ws.Range("D17").Value = "1.935.85"
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").Value = "36.443.98"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.01%  "
$ws.Range("D20").Value = "0.0₃0862"
$ws.Range("E20").Value = "  -3.02%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.79%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.33%  "
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.13"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.134"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.22"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.14%  "
$ws.Range("E30").Value = "  -1.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -6.74%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.14"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -7.04%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.80%  "
$ws.Range("E38").Value = "  -2.56%  "
$ws.Range("E39").Value = "  +4.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0992"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("E42").Value = "  -2.28%  "
$ws.Range("E43").Value = "  -5.87%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.55"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.54%  "
$ws.Range("D45").Value = "1.338.00"
$ws.Range("E45").Value = "  -2.24%  "
$ws.Range("E46").Value = "  -6.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.98"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.65%  "
$ws.Range("E49").Value = "  -0.03%  "
$ws.Range("D50").Value = "2.114.97"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.26"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.73%  "
